# Generate Report for Handoff
# Adds two new handed-off files (c347f44c... and e613a103...) as new rows
# on the Overview / zh-cn / de-de sheets, extending the tables accordingly.

$wb = $excel.ActiveWorkbook

$sheetOverview = $wb.Worksheets.Item("Overview")
$sheetZhCn     = $wb.Worksheets.Item("zh-cn")
$sheetDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet (sheet1) -- columns A:G, new rows 6 & 7
# ---------------------------------------------------------------------------
$ws = $sheetOverview

$ws.Range("A6").Value = "c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md"
$ws.Range("B6").Value = "e2e\c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md"
$ws.Range("B6").Style = "HyperLink"
$ws.Range("C6").Value = ".md"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "Ready for handoff"
$ws.Range("F6").Value = "Ready for handoff"
$ws.Range("G6").Value = "2016-08-24 14:46:39"
$ws.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Range("A7").Value = "e613a103-cebc-48aa-8ef2-2d8405b60fe0.md"
$ws.Range("B7").Value = "e2e\e613a103-cebc-48aa-8ef2-2d8405b60fe0.md"
$ws.Range("B7").Style = "HyperLink"
$ws.Range("C7").Value = ".md"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = "Ready for handoff"
$ws.Range("F7").Value = "Ready for handoff"
$ws.Range("G7").Value = "2016-08-24 14:46:39"
$ws.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws.Hyperlinks.Add($ws.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c347f44c6bb94e38ae798e849f9ea4f8/e2e/c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md", "", "", "e2e\c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e613a103cebc48aa8ef22d8405b60fe0/e2e/e613a103-cebc-48aa-8ef2-2d8405b60fe0.md", "", "", "e2e\e613a103-cebc-48aa-8ef2-2d8405b60fe0.md") | Out-Null

$loOverview = $ws.ListObjects.Item(1)
$loOverview.Resize($ws.Range("A1:G7"))

# ---------------------------------------------------------------------------
# zh-cn sheet (sheet2) -- columns A:P, new rows 6 & 7
# ---------------------------------------------------------------------------
$ws = $sheetZhCn

$ws.Range("A6").Value = "c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md"
$ws.Range("A6").Style = "HyperLink"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "e2e"
$ws.Range("E6").Value = "ht"
$ws.Range("F6").Value = "False"
$ws.Range("G6").Value = "c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.6c9e439a9938fdca1d92c5f02bb0c087af910c71.zh-cn.xlf"
$ws.Range("H6").Value = "2016-08-24 14:46:34"
$ws.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = "0001-01-01 00:00:00"
$ws.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = "True"
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = "False"
$ws.Range("P6").Value = ""

$ws.Range("A7").Value = "e613a103-cebc-48aa-8ef2-2d8405b60fe0.md"
$ws.Range("A7").Style = "HyperLink"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "e2e"
$ws.Range("E7").Value = "ht"
$ws.Range("F7").Value = "False"
$ws.Range("G7").Value = "e613a103-cebc-48aa-8ef2-2d8405b60fe0.0f86b56409d396a006972384bf3e050712d74984.zh-cn.xlf"
$ws.Range("H7").Value = "2016-08-24 14:46:34"
$ws.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = "0001-01-01 00:00:00"
$ws.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = "True"
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = "False"
$ws.Range("P7").Value = ""

$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c347f44c6bb94e38ae798e849f9ea4f8/e2e/c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md", "", "", "c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e613a103cebc48aa8ef22d8405b60fe0/e2e/e613a103-cebc-48aa-8ef2-2d8405b60fe0.md", "", "", "e613a103-cebc-48aa-8ef2-2d8405b60fe0.md") | Out-Null

$loZhCn = $ws.ListObjects.Item(1)
$loZhCn.Resize($ws.Range("A1:P7"))

# ---------------------------------------------------------------------------
# de-de sheet (sheet3) -- columns A:P, new rows 6 & 7
# ---------------------------------------------------------------------------
$ws = $sheetDeDe

$ws.Range("A6").Value = "c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md"
$ws.Range("A6").Style = "HyperLink"
$ws.Range("B6").Value = ".md"
$ws.Range("C6").Value = "Ready for handoff"
$ws.Range("D6").Value = "e2e"
$ws.Range("E6").Value = "ht"
$ws.Range("F6").Value = "False"
$ws.Range("G6").Value = "c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.6c9e439a9938fdca1d92c5f02bb0c087af910c71.de-de.xlf"
$ws.Range("H6").Value = "2016-08-24 14:46:39"
$ws.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = "0001-01-01 00:00:00"
$ws.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = "True"
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = "False"
$ws.Range("P6").Value = ""

$ws.Range("A7").Value = "e613a103-cebc-48aa-8ef2-2d8405b60fe0.md"
$ws.Range("A7").Style = "HyperLink"
$ws.Range("B7").Value = ".md"
$ws.Range("C7").Value = "Ready for handoff"
$ws.Range("D7").Value = "e2e"
$ws.Range("E7").Value = "ht"
$ws.Range("F7").Value = "False"
$ws.Range("G7").Value = "e613a103-cebc-48aa-8ef2-2d8405b60fe0.0f86b56409d396a006972384bf3e050712d74984.de-de.xlf"
$ws.Range("H7").Value = "2016-08-24 14:46:39"
$ws.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = "0001-01-01 00:00:00"
$ws.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = "True"
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = "False"
$ws.Range("P7").Value = ""

$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c347f44c6bb94e38ae798e849f9ea4f8/e2e/c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md", "", "", "c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/e613a103cebc48aa8ef22d8405b60fe0/e2e/e613a103-cebc-48aa-8ef2-2d8405b60fe0.md", "", "", "e613a103-cebc-48aa-8ef2-2d8405b60fe0.md") | Out-Null

$loDeDe = $ws.ListObjects.Item(1)
$loDeDe.Resize($ws.Range("A1:P7"))

Write-Output "Report generated for handback: 2 new files added (c347f44c-6bb9-4e38-ae79-8e849f9ea4f8.md, e613a103-cebc-48aa-8ef2-2d8405b60fe0.md)"
